$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-84 down to 12-85.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new data record.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44635
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112021
$ws.Range("G11").Value = "Ají"
$ws.Range("H11").Value = "Inferno"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 18933
$ws.Range("N11").Value = "`$/caja 12 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 1578
$ws.Range("Q11").Value = 12
$ws.Range("R11").Value = "Hortaliza"

# Ensure the date cell keeps the expected date number-format style (same as
# the rest of column D) after the row insert.
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
